$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 7.572976999999999
$ws.Cells.Item(2, 8).Value = 22.718931
$ws.Cells.Item(2, 9).Value = 0.4497670593913077
$ws.Cells.Item(2, 10).Value = 0.4497670593913078
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 27.85346633333333
$ws.Cells.Item(2, 14).Value = 83.560399
$ws.Cells.Item(2, 15).Value = 0.3884083358054969
$ws.Cells.Item(2, 16).Value = 0.3884083358054969
$ws.Cells.Item(2, 17).Value = 210.9336599126077
$ws.Cells.Item(2, 18).Value = 1898.402939213469
$ws.Cells.Item(2, 19).Value = 0.1746932750383099
$ws.Cells.Item(2, 20).Value = 0.1746932750383099

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 7.572976999999999
$ws.Cells.Item(3, 8).Value = 22.718931
$ws.Cells.Item(3, 9).Value = 0.4497670593913077
$ws.Cells.Item(3, 10).Value = 0.4497670593913078
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 0.6933116666666667
$ws.Cells.Item(3, 14).Value = 2.079935
$ws.Cells.Item(3, 15).Value = 0.009668025782567244
$ws.Cells.Item(3, 16).Value = 0.009668025782567244
$ws.Cells.Item(3, 17).Value = 5.250433305498333
$ws.Cells.Item(3, 18).Value = 47.25389974948499
$ws.Cells.Item(3, 19).Value = 0.004348359526344616
$ws.Cells.Item(3, 20).Value = 0.004348359526344617

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 7.572976999999999
$ws.Cells.Item(4, 8).Value = 22.718931
$ws.Cells.Item(4, 9).Value = 0.4497670593913077
$ws.Cells.Item(4, 10).Value = 0.4497670593913078
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 43.16503600000001
$ws.Cells.Item(4, 14).Value = 129.495108
$ws.Cells.Item(4, 15).Value = 0.601923638411936
$ws.Cells.Item(4, 16).Value = 0.6019236384119359
$ws.Cells.Item(4, 17).Value = 326.887824832172
$ws.Cells.Item(4, 18).Value = 2941.990423489548
$ws.Cells.Item(4, 19).Value = 0.2707254248266532
$ws.Cells.Item(4, 20).Value = 0.2707254248266532

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 2.784025666666667
$ws.Cells.Item(5, 8).Value = 8.352077
$ws.Cells.Item(5, 9).Value = 0.1653462089435359
$ws.Cells.Item(5, 10).Value = 0.1653462089435359
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 27.85346633333333
$ws.Cells.Item(5, 14).Value = 83.560399
$ws.Cells.Item(5, 15).Value = 0.3884083358054969
$ws.Cells.Item(5, 16).Value = 0.3884083358054969
$ws.Cells.Item(5, 17).Value = 77.54476517763588
$ws.Cells.Item(5, 18).Value = 697.902886598723
$ws.Cells.Item(5, 19).Value = 0.06422184584750676
$ws.Cells.Item(5, 20).Value = 0.06422184584750676

# Row 6
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 2.784025666666667
$ws.Cells.Item(6, 8).Value = 8.352077
$ws.Cells.Item(6, 9).Value = 0.1653462089435359
$ws.Cells.Item(6, 10).Value = 0.1653462089435359
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 0.6933116666666667
$ws.Cells.Item(6, 14).Value = 2.079935
$ws.Cells.Item(6, 15).Value = 0.009668025782567244
$ws.Cells.Item(6, 16).Value = 0.009668025782567244
$ws.Cells.Item(6, 17).Value = 1.930197474999444
$ws.Cells.Item(6, 18).Value = 17.371777274995
$ws.Cells.Item(6, 19).Value = 0.001598571411115856
$ws.Cells.Item(6, 20).Value = 0.001598571411115856

# Row 7
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 2.784025666666667
$ws.Cells.Item(7, 8).Value = 8.352077
$ws.Cells.Item(7, 9).Value = 0.1653462089435359
$ws.Cells.Item(7, 10).Value = 0.1653462089435359
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 43.16503600000001
$ws.Cells.Item(7, 14).Value = 129.495108
$ws.Cells.Item(7, 15).Value = 0.601923638411936
$ws.Cells.Item(7, 16).Value = 0.6019236384119359
$ws.Cells.Item(7, 17).Value = 120.1725681265907
$ws.Cells.Item(7, 18).Value = 1081.553113139316
$ws.Cells.Item(7, 19).Value = 0.09952579168491334
$ws.Cells.Item(7, 20).Value = 0.09952579168491332

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 6.480551000000001
$ws.Cells.Item(8, 8).Value = 19.441653
$ws.Cells.Item(8, 9).Value = 0.3848867316651562
$ws.Cells.Item(8, 10).Value = 0.3848867316651562
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 27.85346633333333
$ws.Cells.Item(8, 14).Value = 83.560399
$ws.Cells.Item(8, 15).Value = 0.3884083358054969
$ws.Cells.Item(8, 16).Value = 0.3884083358054969
$ws.Cells.Item(8, 17).Value = 180.5058090999497
$ws.Cells.Item(8, 18).Value = 1624.552281899547
$ws.Cells.Item(8, 19).Value = 0.1494932149196802
$ws.Cells.Item(8, 20).Value = 0.1494932149196802

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 6.480551000000001
$ws.Cells.Item(9, 8).Value = 19.441653
$ws.Cells.Item(9, 9).Value = 0.3848867316651562
$ws.Cells.Item(9, 10).Value = 0.3848867316651562
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 0.6933116666666667
$ws.Cells.Item(9, 14).Value = 2.079935
$ws.Cells.Item(9, 15).Value = 0.009668025782567244
$ws.Cells.Item(9, 16).Value = 0.009668025782567244
$ws.Cells.Item(9, 17).Value = 4.493041614728334
$ws.Cells.Item(9, 18).Value = 40.437374532555
$ws.Cells.Item(9, 19).Value = 0.00372109484510677
$ws.Cells.Item(9, 20).Value = 0.003721094845106771

# Row 10
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 6.480551000000001
$ws.Cells.Item(10, 8).Value = 19.441653
$ws.Cells.Item(10, 9).Value = 0.3848867316651562
$ws.Cells.Item(10, 10).Value = 0.3848867316651562
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 43.16503600000001
$ws.Cells.Item(10, 14).Value = 129.495108
$ws.Cells.Item(10, 15).Value = 0.601923638411936
$ws.Cells.Item(10, 16).Value = 0.6019236384119359
$ws.Cells.Item(10, 17).Value = 279.7332172148361
$ws.Cells.Item(10, 18).Value = 2517.598954933525
$ws.Cells.Item(10, 19).Value = 0.2316724219003693
$ws.Cells.Item(10, 20).Value = 0.2316724219003693

Write-Output "Updated rows 2-10 with new natmi LR-pair values"
